# Added more header related tests.
# - Row 11 ("Create task") becomes "Create a task" with a new endpoint/JSON body.
# - A new row 14 ("Get all tasks") is added, formatted like row 12.
# - Selection moves to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API_Template")
$ws.Activate()

# --- Row 11: "Create task" -> "Create a task" with a new endpoint & JSON body ---
$ws.Range("A11").Value = "Create a task"
$ws.Range("B11").Value = "https://intelliapi-mockserver.herokuapp.com/tasks"

$newJson = "{`n    ""status"": [`n        ""#status""`n    ],`n    ""name"": ""#name"",`n    ""category"": ""#category"",`n    ""isDeleted"": #isDeleted,`n    ""__v"": #version`n}"
$ws.Range("E11").Value = $newJson

# Row height for row 11 grows to fit the longer JSON payload
$ws.Rows.Item(11).RowHeight = 135

# --- New row 14: "Get all tasks" (same look & feel as row 12) ---
$ws.Range("A12:E12").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A14").Value = "Get all tasks"
$ws.Range("B14").Value = "https://intelliapi-mockserver.herokuapp.com/tasks"
$ws.Range("C14").Value = "GET"
$ws.Rows.Item(14).RowHeight = 15

# --- Rebuild the hyperlinks collection (B11 now points at the new tasks URL,
#     and B14 gets a brand new hyperlink) ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://intelliapi-mockserver.herokuapp.com/auth") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://accounts.google.com/o/oauth2/token") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://v2.convertapi.com/upload") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://gorest.co.in/public-api/users") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://gorest.co.in/public-api/users") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6:B8"), "https://gorest.co.in/public-api/users", $null, $null, "https://gorest.co.in/public-api/users") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://maxsoft-mock-server-demo.web.app/photos") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "https://api.apis.guru/", "version/#jsonFile") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://intelliapi-mockserver.herokuapp.com/tasks") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "https://intelliapi-mockserver.herokuapp.com/tasks") | Out-Null

# --- Selection moves to B17 ---
$ws.Range("B17").Select()
